$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel-compatible ROUND(value, 1) implemented on the full-precision decimal
# expansion of the IEEE-754 double (the engine's built-in [Math]::Round /
# WorksheetFunction.Round round the ~15-digit *display* string, which makes
# true halves like 37.05 - whose nearest double is 37.049999999999997 - round
# the wrong way). Working from the 17-significant-digit round-trip string
# reproduces genuine Excel ROUND() semantics (round-half-away-from-zero on
# the real value).
function RoundTo1Decimal($val) {
    if ($val -lt 0) {
        $neg = $true
        $val = -$val
    } else {
        $neg = $false
    }
    $s = $val.ToString("G17")
    if ($s.Contains("E")) {
        $r = [Math]::Round($val, 1)
        if ($neg) { return -$r } else { return $r }
    }
    $dotIdx = $s.IndexOf(".")
    if ($dotIdx -lt 0) {
        if ($neg) { return -([double]$s) } else { return [double]$s }
    }
    $intPart = $s.Substring(0, $dotIdx)
    $fracPart = $s.Substring($dotIdx + 1)
    while ($fracPart.Length -lt 2) { $fracPart = $fracPart + "0" }
    $keep = $fracPart.Substring(0, 1)
    $nextDigit = [int]([string]$fracPart.Substring(1, 1))
    $intKeep = [int64]$intPart
    $keepDigit = [int]$keep
    if ($nextDigit -ge 5) {
        $keepDigit += 1
        if ($keepDigit -eq 10) {
            $keepDigit = 0
            $intKeep += 1
        }
    }
    $result = $intKeep + ($keepDigit / 10.0)
    if ($neg) { return -$result } else { return $result }
}

# Determine the last used row on the sheet (data occupies rows 1..94)
$lastRow = $ws.UsedRange.Rows.Count

# New column is O (15) - "% aplazado" = ROUND(K/M*100, 1)
$headerCol = 15
$kCol = 11
$mCol = 13

# Copy the header cell N1 formatting (bold, centered, bordered) into O1, then
# set its own text.
$ws.Cells.Item(1, 14).Copy($ws.Cells.Item(1, $headerCol))
$ws.Cells.Item(1, $headerCol).Value = "% aplazado"

for ($r = 2; $r -le $lastRow; $r++) {
    $k = $ws.Cells.Item($r, $kCol).Value2
    $m = $ws.Cells.Item($r, $mCol).Value2

    if ($m -eq $null -or $m -eq "") {
        # Rows without an M value (blank/total rows) get a blank cell,
        # matching the pattern of the other empty cells in that row (L/M/N
        # are blank too). Setting .Value = "" on a never-used cell is a
        # no-op in this engine, so copy the already-blank N cell from the
        # same row to materialize an actual (empty) cell at O.
        $ws.Cells.Item($r, 14).Copy($ws.Cells.Item($r, $headerCol))
    } else {
        $pct = RoundTo1Decimal(($k / $m) * 100)
        $ws.Cells.Item($r, $headerCol).Value = $pct
    }
}
